# "updates for spring 2016"
#
# A new row is inserted right under the header row of Sheet1 to record the
# upcoming "Spring_2016" semester (status "future", priority 0.05). All of
# the previously existing data rows shift down by one row as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2 (existing rows 2-21 shift down to 3-22).
[void]$ws.Rows.Item(2).Insert(-4121)

# Copy the formatting of the row directly below (the row that used to be
# row 2, now row 3) into the freshly inserted row so it matches the rest
# of the data rows' look (fonts/number formats).
[void]$ws.Range("A3:D3").Copy()
[void]$ws.Range("A2:D2").PasteSpecial(-4122)

# Populate the new row with the Spring 2016 entry.
$ws.Range("A2").Value = "Spring_2016"
$ws.Range("B2").Value = "future"
$ws.Range("D2").Value = 0.05

# Reflect the cursor/selection ending up on F5 after the edit.
[void]$ws.Range("F5").Select()
